$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 299
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 299
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 897
$ws.Range("N17").Value = -1233
$ws.Range("M17").ClearContents()

$ws.Range("H34").Value = 1499
$ws.Range("I34").Value = 1499
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1499
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1296
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 1499
$ws.Range("I36").Value = 1499
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1499
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -784
$ws.Range("N36").ClearContents()

$ws.Range("H55").Value = 343.55554
$ws.Range("I55").Value = 318.8
$ws.Range("K55").Value = 318.8
$ws.Range("M55").Value = -104.8

$ws.Range("H58").Value = 1173.75
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300

$ws.Range("H70").Value = 3288.7778
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 3942.7144
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 11828.1432
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -12368.1432

$ws.Range("H73").Value = 3288.7778
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 3942.7144
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 11828.1432
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -13700.1432

$ws.Range("H74").Value = 4501.5
$ws.Range("I74").Value = 4501.5
$ws.Range("K74").Value = 4501.5
$ws.Range("M74").Value = -3565.5

$ws.Range("H77").Value = 4501.5
$ws.Range("I77").Value = 4501.5
$ws.Range("K77").Value = 22507.5
$ws.Range("M77").Value = -17827.5

$ws.Range("H80").Value = 7652.143
$ws.Range("J80").Value = 7628.8335
$ws.Range("L80").Value = 22886.5005
$ws.Range("N80").Value = -24882.5005

$ws.Range("H83").Value = 7652.143
$ws.Range("J83").Value = 7628.8335
$ws.Range("L83").Value = 68659.5015
$ws.Range("N83").Value = -78643.5015

$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 1000
$ws.Range("N92").Value = -3496

$ws.Range("H112").Value = 2817.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3446.6924
$ws.Range("I32").Value = 2987.5652
$ws.Range("K32").Value = 2987.5652
$ws.Range("M32").Value = -2700.5652

$ws.Range("H45").Value = 4837.7144
$ws.Range("J45").Value = 4837.7144
$ws.Range("L45").Value = 4837.7144
$ws.Range("N45").Value = -5591.7144

$ws.Range("H97").Value = 350
$ws.Range("I97").Value = 350
$ws.Range("K97").Value = 350
$ws.Range("M97").Value = 146

$ws.Range("H132").Value = 3913.8
$ws.Range("I132").Value = 3896.5
$ws.Range("K132").Value = 11689.5
$ws.Range("M132").Value = -9159.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 8085.4
$ws.Range("I75").Value = 8085.4
$ws.Range("K75").Value = 8085.4
$ws.Range("M75").Value = -7149.4

$ws.Range("H78").Value = 8085.4
$ws.Range("I78").Value = 8085.4
$ws.Range("K78").Value = 24256.2
$ws.Range("M78").Value = -19576.2

$ws.Range("H86").Value = 4026.25
$ws.Range("I86").Value = 4032.6667
$ws.Range("J86").Value = 4007
$ws.Range("K86").Value = 4032.6667
$ws.Range("L86").Value = 4007
$ws.Range("M86").Value = -2909.6667
$ws.Range("N86").Value = -6253

$ws.Range("H89").Value = 4026.25
$ws.Range("I89").Value = 4032.6667
$ws.Range("J89").Value = 4007
$ws.Range("K89").Value = 20163.3335
$ws.Range("L89").Value = 20035
$ws.Range("M89").Value = -14547.3335
$ws.Range("N89").Value = -31267

$ws.Range("H94").Value = 5000
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902

$ws.Range("H105").Value = 8500
$ws.Range("I105").Value = 7000
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 7000
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -5253
$ws.Range("N105").Value = -13494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1407
$ws.Range("I31").Value = 977.73334
$ws.Range("J31").Value = 2211.875
$ws.Range("K31").Value = 977.73334
$ws.Range("L31").Value = 2211.875
$ws.Range("M31").Value = -682.73334
$ws.Range("N31").Value = -2801.875

$ws.Range("H34").Value = 1407
$ws.Range("I34").Value = 977.73334
$ws.Range("J34").Value = 2211.875
$ws.Range("K34").Value = 977.73334
$ws.Range("L34").Value = 2211.875
$ws.Range("M34").Value = -775.73334
$ws.Range("N34").Value = -2615.875

$ws.Range("H59").Value = 12000
$ws.Range("I59").Value = 12000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 12000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -10855
$ws.Range("N59").ClearContents()

$ws.Range("H134").Value = 1362.4
$ws.Range("I134").Value = 1269.3334
$ws.Range("K134").Value = 3808.0002
$ws.Range("M134").Value = -1273.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 3498.2856
$ws.Range("I60").Value = 3248
$ws.Range("J60").Value = 5000
$ws.Range("K60").Value = 9744
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = -9493
$ws.Range("N60").Value = -15502

$ws.Range("H68").Value = 3492.5454
$ws.Range("I68").Value = 4275.143
$ws.Range("J68").Value = 2123
$ws.Range("K68").Value = 12825.429
$ws.Range("L68").Value = 6369
$ws.Range("M68").Value = -12014.429
$ws.Range("N68").Value = -7991

$ws.Range("H71").Value = 3492.5454
$ws.Range("I71").Value = 4275.143
$ws.Range("J71").Value = 2123
$ws.Range("K71").Value = 38476.287
$ws.Range("L71").Value = 19107
$ws.Range("M71").Value = -34420.287
$ws.Range("N71").Value = -27219

$ws.Range("H102").Value = 9999
$ws.Range("I102").Value = 9999
$ws.Range("K102").Value = 29997
$ws.Range("M102").Value = -27563

$ws.Range("H107").Value = 1131.5
$ws.Range("J107").Value = 999.6667
$ws.Range("L107").Value = 2999.0001
$ws.Range("N107").Value = -6839.0001

$ws.Range("H129").Value = 1456.25
$ws.Range("J129").Value = 2616.6667
$ws.Range("L129").Value = 7850.000100000001
$ws.Range("N129").Value = -17850.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 24666.666
$ws.Range("I15").Value = 24000
$ws.Range("K15").Value = 24000
$ws.Range("M15").Value = -23712

$ws.Range("H81").Value = 24666.666
$ws.Range("I81").Value = 24000
$ws.Range("K81").Value = 24000
$ws.Range("M81").Value = -23002

$ws.Range("H84").Value = 24666.666
$ws.Range("I84").Value = 24000
$ws.Range("K84").Value = 72000
$ws.Range("M84").Value = -67008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3365.05
$ws.Range("I46").Value = 2233.4443
$ws.Range("K46").Value = 2233.4443
$ws.Range("M46").Value = -2045.4443

$ws.Range("H61").Value = 2049
$ws.Range("I61").Value = 2049
$ws.Range("K61").Value = 2049
$ws.Range("M61").Value = -1847

$ws.Range("H93").Value = 3000
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -5496

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H113").Value = 2049
$ws.Range("I113").Value = 2049
$ws.Range("K113").Value = 2049
$ws.Range("M113").Value = 121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3857.1428
$ws.Range("J13").Value = 3857.1428
$ws.Range("L13").Value = 3857.1428
$ws.Range("N13").Value = -4137.1428

$ws.Range("H40").Value = 2999.5
$ws.Range("I40").Value = 2999
$ws.Range("K40").Value = 2999
$ws.Range("M40").Value = -2850

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H81").Value = 3340
$ws.Range("I81").Value = 3340
$ws.Range("K81").Value = 6680
$ws.Range("M81").Value = -5619

$ws.Range("H84").Value = 3340
$ws.Range("I84").Value = 3340
$ws.Range("K84").Value = 33400
$ws.Range("M84").Value = -28096
